$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.215.75"
$ws.Range("E2").Value = "  -0.77%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.433.22"
$ws.Range("E3").Value = "  -3.05%  "

# Row 4
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.65"
$ws.Range("E5").Value = "  -3.01%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.23"
$ws.Range("E6").Value = "  -4.94%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.433.34"
$ws.Range("E7").Value = "  -3.07%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.493"
$ws.Range("E9").Value = "  -3.77%  "

# Row 10
$ws.Range("E10").Value = "  -9.56%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.09"
$ws.Range("E11").Value = "  -10.21%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.375"
$ws.Range("E12").Value = "  -7.18%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.015.21"
$ws.Range("E13").Value = "  -3.05%  "

# Row 14
$ws.Range("E14").Value = "  -9.34%  "

# Row 15
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.49"
$ws.Range("E15").Value = "  -7.13%  "

# Row 16
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.439.40"
$ws.Range("E16").Value = "  -3.40%  "

# Row 17
$ws.Range("E17").Value = "  -1.96%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "65.222.67"
$ws.Range("E18").Value = "  -0.70%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.65"
$ws.Range("E19").Value = "  -12.45%  "

# Row 20
$ws.Range("E20").Value = "  -6.10%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.48"
$ws.Range("E21").Value = "  -5.85%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "385.09"
$ws.Range("E22").Value = "  -7.87%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.552"
$ws.Range("E23").Value = "  -7.83%  "

# Row 24
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  -0.12%  "

# Row 25
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.54"
$ws.Range("E25").Value = "  -6.80%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.572.02"
$ws.Range("E26").Value = "  -2.93%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000106"
$ws.Range("E27").Value = "  -8.20%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  +0.01%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.15"
$ws.Range("E29").Value = "  -8.34%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.09"
$ws.Range("E30").Value = "  -8.88%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.20"
$ws.Range("E31").Value = "  -10.39%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.441.77"
$ws.Range("E32").Value = "  -3.06%  "

# Row 33
$ws.Range("E33").Value = "  +0.01%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.144"
$ws.Range("E34").Value = "  -6.97%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.07"
$ws.Range("E35").Value = "  -5.23%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "168.97"
$ws.Range("E36").Value = "  -3.10%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.76"
$ws.Range("E37").Value = "  -10.46%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.17"
$ws.Range("E38").Value = "  -11.28%  "

# Row 39
$ws.Range("E39").Value = "  -7.43%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.67"
$ws.Range("E40").Value = "  -11.38%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0755"
$ws.Range("E41").Value = "  -7.52%  "

# Row 42
$ws.Range("E42").Value = "  -5.47%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.66"
$ws.Range("E43").Value = "  -5.43%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.15%  "

# Row 45
$ws.Range("E45").Value = "  -14.32%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.61"
$ws.Range("E46").Value = "  -9.28%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.82"
$ws.Range("E48").Value = "  -2.32%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.47"
$ws.Range("E49").Value = "  -7.78%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.04"
$ws.Range("E50").Value = "  -13.43%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.164.89"
$ws.Range("E51").Value = "  -7.30%  "
